# Swap the order of slide 2 ("What is Physical Computing?") and
# slide 3 ("Interactive..." diagram) so the diagram slide now comes
# second and the "What is Physical Computing?" slide comes third.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.MoveTo(3)
